$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 859; this shifts every existing
# row from 859..900 down to 860..901 and expands the used range /
# dimension automatically (A1:D900 -> A1:D901).
$ws.Rows("859:859").Insert()

# Fill in the new row's data: 2026/02/26 (Thu), hour 2, ranking 22.
# Column A holds the date as literal text (matching every other date
# cell in the sheet), so force a Text number format before assigning
# the value to stop Excel from auto-converting the "yyyy/mm/dd"
# looking string into a real date serial number. Reset the format
# back to Normal afterwards so the cell doesn't end up with a
# different style than its neighbours.
$ws.Range("A859").NumberFormat = "@"
$ws.Range("A859").Value = "2026/02/26"
$ws.Range("A859").Style = "Normal"

$ws.Range("B859").Value = "木"
$ws.Range("C859").Value = 2
$ws.Range("D859").Value = 22
